$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update J2:J6 from 32 -> 8
$ws.Range("J2:J6").Value = 8

# Update K2:K6 from 16 -> 4
$ws.Range("K2:K6").Value = 4

# Update selection to K2:K6 with active cell K2
$ws.Range("K2:K6").Select()

# Update window position (xWindow/yWindow)
$win = $excel.ActiveWindow
$win.Left = 5760
$win.Top = 1740
